# Update gh-pages output (丽水 conventions data) for sheets "展览" and "全部类型".
# Rows 2-5 get new event data, row 6 is removed entirely (5 events -> 4 events).

$wb = $excel.ActiveWorkbook

$targetSheets = @("展览", "全部类型")

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 2
    # Leading apostrophe forces text so Excel doesn't auto-convert the
    # yyyy-mm-dd-looking string into a real date serial number.
    $ws.Range("B2").Value = "'2024-02-07"
    $ws.Range("C2").Value = "丽水·新年动漫狂欢盛典"
    $ws.Range("D2").Value = "中东路848号(解放街交汇) 飞达国际大酒店"
    $ws.Range("E2").Value = "2024.02.07 09:00-02.07 17:00"
    $ws.Range("F2").Value = 301
    $ws.Range("G2").Value = 45
    $ws.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=78294"
    $ws.Range("I2").Value = "//i1.hdslb.com/bfs/openplatform/202311/lP5IkqWn1699431829470.jpeg"

    # Row 3
    $ws.Range("B3").Value = "'2024-02-07"
    $ws.Range("C3").Value = "龙泉·崩X铁X原ONLY"
    $ws.Range("D3").Value = "金沙路26-1号 龙泉金沙温泉酒店"
    $ws.Range("E3").Value = "2024.02.07 10:30-02.07 16:30"
    $ws.Range("F3").Value = 234
    $ws.Range("G3").Value = 50
    $ws.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=80714"
    $ws.Range("I3").Value = "//i2.hdslb.com/bfs/openplatform/202401/rTvQio211704877379770.jpeg"

    # Row 4
    $ws.Range("B4").Value = "'2024-02-14"
    $ws.Range("C4").Value = "丽水·YA●怀旧only"
    $ws.Range("D4").Value = "人民街567号 丽水体育中心"
    $ws.Range("E4").Value = "2024.02.14 09:00-02.14 17:00"
    $ws.Range("F4").Value = 43
    $ws.Range("G4").Value = 45
    $ws.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=81032"
    $ws.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202401/LbqTNkvq1705561884633.png"

    # Row 5
    $ws.Range("B5").Value = "'2024-02-18"
    $ws.Range("C5").Value = "丽水·LPJ 现实X次元动漫展"
    $ws.Range("D5").Value = "中东路848号(解放街交汇) 飞达国际大酒店"
    $ws.Range("E5").Value = "2024.02.18 10:00-02.18 17:00"
    $ws.Range("F5").Value = 270
    $ws.Range("G5").Value = 45
    $ws.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=79437"
    $ws.Range("I5").Value = "//i1.hdslb.com/bfs/openplatform/202312/ee5hLUN61702276208812.jpeg"

    # Row 6 (缙云·星辰动漫游戏展嘉年华) no longer exists upstream - remove the row entirely.
    $ws.Rows("6:6").Delete()
}
